$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Status" header in column F (new employee status filter column)
$ws.Range("F1").Value = "Status"

# Update the active selection to the newly added cell
$ws.Range("F1").Select()
